# Updates cryptos list price (D) and 1h volume % (E) columns for rows 2-51.
# D-column values are plain text (e.g. "593.93"); many of them are
# syntactically valid numbers, so a bare Range.Value assignment would get
# auto-coerced to a float (losing the text type / trailing zeros). We force
# text entry with NumberFormat="@" and then reset the style back to "Normal"
# afterwards so no stray style index is left on the cell (matches the
# original file, where these cells carry no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.412.83'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.771.90'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.769.96'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  -0.16%  '
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.39'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.96%  '
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("E13").Value = '  -2.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.404.86'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.747.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.487.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.76%  '
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '457.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.47%  '
$ws.Range("E23").Value = '  -0.87%  '
$ws.Range("E24").Value = '  +6.23%  '
$ws.Range("E25").Value = '  -1.77%  '
$ws.Range("E26").Value = '  -4.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("E28").Value = '  -1.70%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.67%  '
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.727.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  -1.58%  '
$ws.Range("E38").Value = '  -1.23%  '
$ws.Range("E39").Value = '  -1.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.993'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '44.10'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("E45").Value = '  -2.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.83'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.55%  '
$ws.Range("E47").Value = '  -3.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '145.92'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '392.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.755.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.80%  '
